$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1120.4615
$ws.Range("J17").Value = 1132.5217
$ws.Range("L17").Value = 3397.5651
$ws.Range("N17").Value = -3733.5651

$ws.Range("H41").Value = 395.41666
$ws.Range("I41").Value = 202.875
$ws.Range("J41").Value = 780.5
$ws.Range("K41").Value = 202.875
$ws.Range("L41").Value = 780.5
$ws.Range("M41").Value = 237.125
$ws.Range("N41").Value = -1660.5

$ws.Range("H74").Value = 6912.375
$ws.Range("I74").Value = 3766.3333
$ws.Range("K74").Value = 3766.3333
$ws.Range("M74").Value = -2830.3333

$ws.Range("H77").Value = 6912.375
$ws.Range("I77").Value = 3766.3333
$ws.Range("K77").Value = 18831.6665
$ws.Range("M77").Value = -14151.6665

$ws.Range("H98").Value = 17382.17
$ws.Range("I98").Value = 18240.387
$ws.Range("J98").Value = 14721.7
$ws.Range("K98").Value = 18240.387
$ws.Range("L98").Value = 14721.7
$ws.Range("M98").Value = -16742.387
$ws.Range("N98").Value = -17717.7

$ws.Range("H122").Value = 17382.17
$ws.Range("I122").Value = 18240.387
$ws.Range("J122").Value = 14721.7
$ws.Range("K122").Value = 54721.16099999999
$ws.Range("L122").Value = 44165.10000000001
$ws.Range("M122").Value = -52271.16099999999
$ws.Range("N122").Value = -49065.10000000001

$ws.Range("H133").Value = 98203
$ws.Range("J133").Value = 98203
$ws.Range("L133").Value = 98203
$ws.Range("N133").Value = -108323

$ws.Range("H138").Value = 3903.1968
$ws.Range("I138").Value = 1426.1177
$ws.Range("J138").Value = 4860.25
$ws.Range("K138").Value = 4278.3531
$ws.Range("L138").Value = 14580.75
$ws.Range("M138").Value = 861.6468999999997
$ws.Range("N138").Value = -24860.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2974.831
$ws.Range("I32").Value = 2589.9177
$ws.Range("K32").Value = 2589.9177
$ws.Range("M32").Value = -2302.9177

$ws.Range("H63").Value = 4558.9
$ws.Range("I63").Value = 3234.1428
$ws.Range("K63").Value = 3234.1428
$ws.Range("M63").Value = -2548.1428

$ws.Range("H66").Value = 4558.9
$ws.Range("I66").Value = 3234.1428
$ws.Range("K66").Value = 16170.714
$ws.Range("M66").Value = -12738.714

$ws.Range("H74").Value = 12122
$ws.Range("I74").Value = 1110
$ws.Range("J74").Value = 14875
$ws.Range("K74").Value = 1110
$ws.Range("L74").Value = 14875
$ws.Range("M74").Value = -236
$ws.Range("N74").Value = -16623

$ws.Range("H75").Value = 55000
$ws.Range("J75").Value = 55000
$ws.Range("L75").Value = 55000
$ws.Range("N75").Value = -56748

$ws.Range("H77").Value = 12122
$ws.Range("I77").Value = 1110
$ws.Range("J77").Value = 14875
$ws.Range("K77").Value = 5550
$ws.Range("L77").Value = 74375
$ws.Range("M77").Value = -1182
$ws.Range("N77").Value = -83111

$ws.Range("H78").Value = 55000
$ws.Range("J78").Value = 55000
$ws.Range("L78").Value = 165000
$ws.Range("N78").Value = -173736

$ws.Range("H132").Value = 4145.2236
$ws.Range("I132").Value = 4076.7715
$ws.Range("K132").Value = 12230.3145
$ws.Range("M132").Value = -9700.3145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 111214
$ws.Range("I75").Value = 111214
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 111214
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -110278
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 111214
$ws.Range("I78").Value = 111214
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 333642
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -328962
$ws.Range("N78").ClearContents()

$ws.Range("H94").Value = 3489.0557
$ws.Range("I94").Value = 2349.8333
$ws.Range("J94").Value = 5767.5
$ws.Range("K94").Value = 2349.8333
$ws.Range("L94").Value = 5767.5
$ws.Range("M94").Value = -1898.8333
$ws.Range("N94").Value = -6669.5

$ws.Range("H118").Value = 8400
$ws.Range("J118").Value = 8400
$ws.Range("L118").Value = 8400
$ws.Range("N118").Value = -11714

$ws.Range("H134").Value = 18756.938
$ws.Range("I134").Value = 23050.916
$ws.Range("K134").Value = 69152.74800000001
$ws.Range("M134").Value = -66617.74800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4896.9165
$ws.Range("I31").Value = 4951.4443
$ws.Range("K31").Value = 4951.4443
$ws.Range("M31").Value = -4656.4443

$ws.Range("H34").Value = 4896.9165
$ws.Range("I34").Value = 4951.4443
$ws.Range("K34").Value = 4951.4443
$ws.Range("M34").Value = -4749.4443

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H86").Value = 6780.125
$ws.Range("I86").Value = 5129.385
$ws.Range("K86").Value = 5129.385
$ws.Range("M86").Value = -4006.385

$ws.Range("H89").Value = 6780.125
$ws.Range("I89").Value = 5129.385
$ws.Range("K89").Value = 25646.925
$ws.Range("M89").Value = -20030.925

$ws.Range("H97").Value = 63000
$ws.Range("J97").Value = 66000
$ws.Range("L97").Value = 66000
$ws.Range("N97").Value = -67982

$ws.Range("H121").Value = 55000
$ws.Range("J121").Value = 55000
$ws.Range("L121").Value = 55000
$ws.Range("N121").Value = -57620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 325140.7
$ws.Range("J5").Value = 387340.84
$ws.Range("L5").Value = 1162022.52
$ws.Range("N5").Value = -1162246.52

$ws.Range("H55").Value = 10051.853
$ws.Range("J55").Value = 10570.594
$ws.Range("L55").Value = 31711.782
$ws.Range("N55").Value = -32065.782

$ws.Range("H107").Value = 562.1111
$ws.Range("I107").Value = 303.84616
$ws.Range("J107").Value = 644
$ws.Range("K107").Value = 911.5384799999999
$ws.Range("L107").Value = 1932
$ws.Range("M107").Value = 1008.46152
$ws.Range("N107").Value = -5772

$ws.Range("H135").Value = 325140.7
$ws.Range("J135").Value = 387340.84
$ws.Range("L135").Value = 3486067.56
$ws.Range("N135").Value = -3491137.56

$ws.Range("H140").Value = 353136
$ws.Range("I140").Value = 353136
$ws.Range("K140").Value = 1059408
$ws.Range("M140").Value = -1054228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 5406779.5
$ws.Range("I18").Value = 9008666
$ws.Range("J18").Value = 3949.5
$ws.Range("K18").Value = 9008666
$ws.Range("L18").Value = 3949.5
$ws.Range("M18").Value = -9008373
$ws.Range("N18").Value = -4535.5

$ws.Range("H68").Value = 30000
$ws.Range("I68").Value = 30000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 30000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -29189
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 30000
$ws.Range("I71").Value = 30000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 90000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -85944
$ws.Range("N71").ClearContents()

$ws.Range("H80").Value = 12155
$ws.Range("J80").Value = 3414.3333
$ws.Range("L80").Value = 3414.3333
$ws.Range("N80").Value = -5410.3333

$ws.Range("H83").Value = 12155
$ws.Range("J83").Value = 3414.3333
$ws.Range("L83").Value = 17071.6665
$ws.Range("N83").Value = -27055.6665

$ws.Range("H88").Value = 25000
$ws.Range("J88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("N88").Value = -25902

$ws.Range("H91").Value = 25000
$ws.Range("J91").Value = 25000
$ws.Range("L91").Value = 25000
$ws.Range("N91").Value = -28120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 75098
$ws.Range("I74").Value = 70196
$ws.Range("K74").Value = 70196
$ws.Range("M74").Value = -69198

$ws.Range("H77").Value = 75098
$ws.Range("I77").Value = 70196
$ws.Range("K77").Value = 210588
$ws.Range("M77").Value = -205596

$ws.Range("H136").Value = 3966.1714
$ws.Range("I136").Value = 1311.2413
$ws.Range("J136").Value = 16798.334
$ws.Range("K136").Value = 3933.7239
$ws.Range("L136").Value = 50395.00199999999
$ws.Range("M136").Value = -1383.7239
$ws.Range("N136").Value = -55495.00199999999

$ws.Range("H137").Value = 75000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 75000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 168006.48
$ws.Range("I62").Value = 635524.8
$ws.Range("K62").Value = 635524.8
$ws.Range("M62").Value = -634900.8

$ws.Range("H65").Value = 168006.48
$ws.Range("I65").Value = 635524.8
$ws.Range("K65").Value = 3177624
$ws.Range("M65").Value = -3174504

$ws.Range("H132").Value = 13483.808
$ws.Range("I132").Value = 14789.714
$ws.Range("K132").Value = 44369.142
$ws.Range("M132").Value = -41839.142

$ws.Range("H136").Value = 2051.9143
$ws.Range("I136").Value = 1390.1072
$ws.Range("J136").Value = 4699.143
$ws.Range("K136").Value = 4170.321599999999
$ws.Range("L136").Value = 14097.429
$ws.Range("M136").Value = -1620.321599999999
$ws.Range("N136").Value = -19197.429

